# Add team win/loss/tie record columns (Wins, Losses, Ties) right after the
# existing data (columns A:AC), as described in the commit message:
# "Added team record to data ... the W/L/T are on the same sheet".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells live in AD1:AF1 (columns 30-32), right after the existing
# "Unnamed: 28" header in AC1 (column 29).
# Copy the formatting of an existing header cell (A1 - bold, centered, top
# aligned, thin border) onto the new header cells so they match the rest of
# row 1, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-38) gets the same team record for the season:
# 53 wins, 62 losses, 0 ties.
for ($r = 2; $r -le 38; $r++) {
    $ws.Cells.Item($r, 30).Value = 53   # AD: Wins
    $ws.Cells.Item($r, 31).Value = 62   # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF: Ties
}
